# The workbook's single worksheet ("컬럼정보" - column info) documents the
# columns of a dataset. Four of the description cells (C5:C8) described the
# "representative floor-usage" columns generically; the author clarified that
# these values exclude basement floors by prefixing the description text
# with "(지하층 제외) " (i.e. "excluding basement floors").

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("C5").Value = "(지하층 제외) 해당 건물에서 면적 합이 가장 큰 층용도 코드"
$ws.Range("C6").Value = "(지하층 제외) 해당 건물에서 면적 합이 가장 큰 층용도 명"
$ws.Range("C7").Value = "(지하층 제외) 해당 건물에서 면적 합이 가장 큰 층용도 면적"
$ws.Range("C8").Value = "(지하층 제외) 해당 건물에서 면적 합이 가장 큰 층용도 면적의 비율"

# Selection in the saved file moved from C10 to C9.
$ws.Range("C9").Select()
